$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the cryptos-list price/volume refresh (GitHub Actions bot update).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# Numeric-looking D values are written via a temporary Text number format
# so COM does not silently coerce them into Number cells (and drop the
# trailing zero), then the format is cleared to match the original
# (unstyled) inlineStr cells.

$ws.Range('D2').Value = '57.904.32'
$ws.Range('E2').Value = '  +0.95%  '

$ws.Range('D3').Value = '3.133.85'
$ws.Range('E3').Value = '  +1.06%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '533.78'
$ws.Range('D5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.50'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.67%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').Value = '3.132.07'
$ws.Range('E8').Value = '  +1.03%  '

$ws.Range('E9').Value = '  +5.68%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.34'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.19%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.108'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.94%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.417'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.26%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.137'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.15%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.661.86'
$ws.Range('E14').Value = '  +0.75%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.97'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.99%  '

$ws.Range('E16').Value = '  +1.56%  '

$ws.Range('D17').Value = '57.987.78'
$ws.Range('E17').Value = '  +0.90%  '

$ws.Range('D18').Value = '3.131.03'
$ws.Range('E18').Value = '  +1.17%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.10'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.26%  '

$ws.Range('E20').Value = '  +3.12%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.11'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.56%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '368.80'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +6.15%  '

$ws.Range('E23').Value = '  -0.22%  '

$ws.Range('E24').Value = '  -2.66%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '69.32'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.70%  '

$ws.Range('E26').Value = '  +1.79%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.169'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.40%  '

$ws.Range('E28').Value = '  -1.73%  '

$ws.Range('D29').Value = '0.0₃0868'
$ws.Range('E29').Value = '  -2.42%  '

$ws.Range('E30').Value = '  -0.57%  '

$ws.Range('E31').Value = '  +0.21%  '

$ws.Range('E32').Value = '  +2.02%  '

$ws.Range('E33').Value = '  +3.86%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.18'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.14%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.19'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.72%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.46'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.64%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.10'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.08%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.30'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.53%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.60'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.39%  '

$ws.Range('E40').Value = '  +3.67%  '

$ws.Range('E41').Value = '  +2.51%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.12'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.86%  '

$ws.Range('D43').Value = '2.516.36'
$ws.Range('E43').Value = '  +5.66%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.701'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.68%  '

$ws.Range('E45').Value = '  +3.83%  '

$ws.Range('E46').Value = '  +2.17%  '

$ws.Range('E47').Value = '  -0.04%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.983'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.75%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.15'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.60%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.78'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.87%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.748'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.10%  '

